{"js": "// New answer values for the 20x5 table of arithmetic problems, in row-major\n// order (document order: row 0 col 0..4, row 1 col 0..4, ...).\nconst newValues = [\n  [\"6-5=1\", \"49-17=32\", \"73+4=77\", \"29+31=60\", \"50+12=62\"],\n  [\"33-27=6\", \"44-31=13\", \"96-75=21\", \"37+21=58\", \"6+33=39\"],\n  [\"52-31=21\", \"92-68=24\", \"52-14=38\", \"29+60=89\", \"65+20=85\"],\n  [\"6+5=11\", \"77-2=75\", \"45+42=87\", \"24+54=78\", \"4+42=46\"],\n  [\"77-55=22\", \"37-26=11\", \"82-76=6\", \"70-26=44\", \"75-30=45\"],\n  [\"82-38=44\", \"78-74=4\", \"50+0=50\", \"54-41=13\", \"68-30=38\"],\n  [\"84-23=61\", \"22+8=30\", \"91-55=36\", \"59+10=69\", \"93-80=13\"],\n  [\"38+10=48\", \"13+46=59\", \"19-2=17\", \"37-24=13\", \"91-78=13\"],\n  [\"72+14=86\", \"32+15=47\", \"3+63=66\", \"75+8=83\", \"43+1=44\"],\n  [\"12+18=30\", \"65-56=9\", \"68+10=78\", \"76+8=84\", \"85-7=78\"],\n  [\"35-0=35\", \"75-49=26\", \"87-9=78\", \"42-13=29\", \"18+15=33\"],\n  [\"5+33=38\", \"93-20=73\", \"66-60=6\", \"80-6=74\", \"49+1=50\"],\n  [\"33-21=12\", \"55+4=59\", \"89-14=75\", \"95-13=82\", \"66+25=91\"],\n  [\"13+6=19\", \"9+29=38\", \"67-55=12\", \"59+34=93\", \"61+9=70\"],\n  [\"55+33=88\", \"16+80=96\", \"39+48=87\", \"5+63=68\", \"75-46=29\"],\n  [\"71+19=90\", \"9+84=93\", \"85+6=91\", \"29+31=60\", \"44-25=19\"],\n  [\"31+52=83\", \"18+46=64\", \"78-3=75\", \"9+46=55\", \"92-66=26\"],\n  [\"99+0=99\", \"15+65=80\", \"38+8=46\", \"1+6=7\", \"37+4=41\"],\n  [\"12+56=68\", \"31+55=86\", \"46-38=8\", \"65-3=62\", \"83-62=21\"],\n  [\"97-32=65\", \"21-4=17\", \"73-28=45\", \"8-7=1\", \"52-30=22\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\n// Replace the text of each cell's (single) paragraph range in place so the\n// existing run formatting (font, size) and paragraph formatting (alignment)\n// are preserved \u2014 replacing at the cell-body level instead would wipe the\n// run/paragraph properties.\nfor (let r = 0; r < newValues.length; r++) {\n  for (let c = 0; c < newValues[r].length; c++) {\n    const cell = table.getCell(r, c);\n    const para = cell.body.paragraphs.getFirst();\n    const rng = para.getRange();\n    rng.insertText(newValues[r][c], Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# New answer values for the 20x5 table of arithmetic problems, in row-major\n# order (document order: row 1 col 1..5, row 2 col 1..5, ...). Word COM\n# collections (Rows/Columns/Cell) are 1-based.\n$newValues = @(\n    @(\"6-5=1\", \"49-17=32\", \"73+4=77\", \"29+31=60\", \"50+12=62\"),\n    @(\"33-27=6\", \"44-31=13\", \"96-75=21\", \"37+21=58\", \"6+33=39\"),\n    @(\"52-31=21\", \"92-68=24\", \"52-14=38\", \"29+60=89\", \"65+20=85\"),\n    @(\"6+5=11\", \"77-2=75\", \"45+42=87\", \"24+54=78\", \"4+42=46\"),\n    @(\"77-55=22\", \"37-26=11\", \"82-76=6\", \"70-26=44\", \"75-30=45\"),\n    @(\"82-38=44\", \"78-74=4\", \"50+0=50\", \"54-41=13\", \"68-30=38\"),\n    @(\"84-23=61\", \"22+8=30\", \"91-55=36\", \"59+10=69\", \"93-80=13\"),\n    @(\"38+10=48\", \"13+46=59\", \"19-2=17\", \"37-24=13\", \"91-78=13\"),\n    @(\"72+14=86\", \"32+15=47\", \"3+63=66\", \"75+8=83\", \"43+1=44\"),\n    @(\"12+18=30\", \"65-56=9\", \"68+10=78\", \"76+8=84\", \"85-7=78\"),\n    @(\"35-0=35\", \"75-49=26\", \"87-9=78\", \"42-13=29\", \"18+15=33\"),\n    @(\"5+33=38\", \"93-20=73\", \"66-60=6\", \"80-6=74\", \"49+1=50\"),\n    @(\"33-21=12\", \"55+4=59\", \"89-14=75\", \"95-13=82\", \"66+25=91\"),\n    @(\"13+6=19\", \"9+29=38\", \"67-55=12\", \"59+34=93\", \"61+9=70\"),\n    @(\"55+33=88\", \"16+80=96\", \"39+48=87\", \"5+63=68\", \"75-46=29\"),\n    @(\"71+19=90\", \"9+84=93\", \"85+6=91\", \"29+31=60\", \"44-25=19\"),\n    @(\"31+52=83\", \"18+46=64\", \"78-3=75\", \"9+46=55\", \"92-66=26\"),\n    @(\"99+0=99\", \"15+65=80\", \"38+8=46\", \"1+6=7\", \"37+4=41\"),\n    @(\"12+56=68\", \"31+55=86\", \"46-38=8\", \"65-3=62\", \"83-62=21\"),\n    @(\"97-32=65\", \"21-4=17\", \"73-28=45\", \"8-7=1\", \"52-30=22\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nfor ($r = 0; $r -lt $newValues.Length; $r++) {\n    $row = $newValues[$r]\n    for ($c = 0; $c -lt $row.Length; $c++) {\n        # Setting Range.Text replaces just the cell's text content while Word\n        # keeps the existing run/paragraph formatting (font, size, alignment)\n        # and automatically preserves the end-of-cell marker.\n        $t.Cell($r + 1, $c + 1).Range.Text = $row[$c]\n    }\n}\n"}
